$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data
$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -2
$ws.Range("F8").Value = -2
$ws.Range("F12").Value = 9
$ws.Range("F14").Value = 13
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = 0

$wb.Save()
